$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.037.84"
Set-TextValue $ws.Range("E2") "  +6.95%  "
Set-TextValue $ws.Range("D3") "1.740.86"
Set-TextValue $ws.Range("E3") "  +5.15%  "
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "228.70"
Set-TextValue $ws.Range("E5") "  +4.29%  "
Set-TextValue $ws.Range("D6") "0.5445"
Set-TextValue $ws.Range("E6") "  +3.79%  "
Set-TextValue $ws.Range("E7") "  -0.18%  "
Set-TextValue $ws.Range("D8") "0.2769"
Set-TextValue $ws.Range("E8") "  +3.78%  "
Set-TextValue $ws.Range("D9") "0.06742"
Set-TextValue $ws.Range("E9") "  +6.02%  "
Set-TextValue $ws.Range("D10") "21.70"
Set-TextValue $ws.Range("E10") "  +4.86%  "
Set-TextValue $ws.Range("D11") "0.07783"
Set-TextValue $ws.Range("E11") "  +0.78%  "
Set-TextValue $ws.Range("D12") "4.703"
Set-TextValue $ws.Range("E12") "  +2.32%  "
Set-TextValue $ws.Range("D13") "1.730.52"
Set-TextValue $ws.Range("E13") "  +8.08%  "
Set-TextValue $ws.Range("D14") "1.981.29"
Set-TextValue $ws.Range("E14") "  +5.21%  "
Set-TextValue $ws.Range("D15") "0.5989"
Set-TextValue $ws.Range("E15") "  +6.24%  "
Set-TextValue $ws.Range("D16") "0.0₅8408"
Set-TextValue $ws.Range("E16") "  +1.78%  "
Set-TextValue $ws.Range("D17") "69.06"
Set-TextValue $ws.Range("E17") "  +5.58%  "
Set-TextValue $ws.Range("D18") "28.006.67"
Set-TextValue $ws.Range("E18") "  +6.84%  "
Set-TextValue $ws.Range("D19") "224.93"
Set-TextValue $ws.Range("D20") "4.842"
Set-TextValue $ws.Range("E20") "  +3.00%  "
Set-TextValue $ws.Range("E21") "  -0.22%  "
Set-TextValue $ws.Range("D22") "10.90"
Set-TextValue $ws.Range("E22") "  +4.85%  "
Set-TextValue $ws.Range("D23") "6.238"
Set-TextValue $ws.Range("E23") "  +3.84%  "
Set-TextValue $ws.Range("E24") "  -0.14%  "
Set-TextValue $ws.Range("D25") "146.19"
Set-TextValue $ws.Range("E25") "  +1.60%  "
Set-TextValue $ws.Range("D26") "0.1251"
Set-TextValue $ws.Range("E26") "  +3.75%  "
Set-TextValue $ws.Range("D27") "7.475"
Set-TextValue $ws.Range("E27") "  +2.79%  "
Set-TextValue $ws.Range("D28") "17.15"
Set-TextValue $ws.Range("E28") "  +7.62%  "
Set-TextValue $ws.Range("D29") "1.646"
Set-TextValue $ws.Range("E29") "  +10.12%  "
Set-TextValue $ws.Range("D30") "0.05696"
Set-TextValue $ws.Range("E30") "  +1.29%  "
Set-TextValue $ws.Range("D31") "1.319"
Set-TextValue $ws.Range("E31") "  +3.43%  "
Set-TextValue $ws.Range("D32") "3.722"
Set-TextValue $ws.Range("E32") "  +6.19%  "
Set-TextValue $ws.Range("D33") "3.530"
Set-TextValue $ws.Range("E33") "  +4.92%  "
Set-TextValue $ws.Range("D34") "1.676"
Set-TextValue $ws.Range("E34") "  +5.79%  "
Set-TextValue $ws.Range("D35") "0.9837"
Set-TextValue $ws.Range("E35") "  +3.28%  "
Set-TextValue $ws.Range("D36") "2.859"
Set-TextValue $ws.Range("D37") "2.452"
Set-TextValue $ws.Range("E37") "  +1.62%  "
Set-TextValue $ws.Range("D38") "0.5965"
Set-TextValue $ws.Range("E38") "  +3.55%  "
Set-TextValue $ws.Range("E39") "  +5.05%  "
Set-TextValue $ws.Range("D40") "5.983"
Set-TextValue $ws.Range("E40") "  -0.65%  "
Set-TextValue $ws.Range("D41") "1.049.82"
Set-TextValue $ws.Range("E41") "  +3.80%  "
Set-TextValue $ws.Range("D42") "0.8482"
Set-TextValue $ws.Range("E42") "  +0.70%  "
Set-TextValue $ws.Range("E43") "  -0.12%  "
Set-TextValue $ws.Range("D44") "102.08"
Set-TextValue $ws.Range("E44") "  +0.10%  "
Set-TextValue $ws.Range("D45") "1.886.47"
Set-TextValue $ws.Range("E45") "  +5.14%  "
Set-TextValue $ws.Range("D46") "0.0₈117"
Set-TextValue $ws.Range("E46") "  +13.47%  "
Set-TextValue $ws.Range("D47") "60.06"
Set-TextValue $ws.Range("E47") "  +2.68%  "
Set-TextValue $ws.Range("D48") "8.307"
Set-TextValue $ws.Range("E48") "  +3.51%  "
Set-TextValue $ws.Range("D49") "0.4429"
Set-TextValue $ws.Range("E49") "  +1.87%  "
Set-TextValue $ws.Range("D50") "1.001"
Set-TextValue $ws.Range("E50") "  -0.77%  "
Set-TextValue $ws.Range("D51") "0.05315"
Set-TextValue $ws.Range("E51") "  -0.54%  "
